$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting from the
# existing header cell H1 (bold, bordered, centered) then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), keyed by row number.
$data = @{
    2  = @(5, 5)
    3  = @(6, 7)
    4  = @(8, 9)
    5  = @(8, 9)
    6  = @(4, 5)
    7  = @(6, 9)
    8  = @(5, 7)
    9  = @(9, 9)
    10 = @(8, 9)
    11 = @(6, 7)
    12 = @(7, 8)
    13 = @(9, 9)
    14 = @(5, 6)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(7, 7)
    19 = @(8, 8)
    20 = @(7, 8)
    21 = @(8, 8)
    22 = @(5, 5)
    23 = @(6, 6)
    24 = @(6, 6)
    25 = @(3, 5)
    26 = @(7, 8)
    27 = @(7, 7)
    28 = @(9, 9)
    29 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
